{"js": "// Update the two-digit multiplication prompts in the worksheet table.\n// Each entry identifies a table cell by its 0-based (row, column) position\n// and gives the expected current text plus the replacement text.\n// Addressing cells by position (rather than a single document-wide\n// search/replace) is required because several prompts share the same\n// original text (e.g. \"94\u00d750=\" and \"74\u00d795=\" each occur twice in the\n// table) but must become different values depending on which cell they\n// are in. Searching/replacing within each individual cell's body keeps\n// the match scoped to exactly one occurrence and preserves the existing\n// run/paragraph formatting (font, size, alignment) because the replaced\n// range sits inside the original run.\n\nconst edits = [\n  { row: 0, col: 0, oldText: \"58\u00d752=\", newText: \"18\u00d779=\" },\n  { row: 0, col: 1, oldText: \"99\u00d787=\", newText: \"46\u00d744=\" },\n  { row: 0, col: 2, oldText: \"77\u00d781=\", newText: \"72\u00d723=\" },\n  { row: 0, col: 3, oldText: \"98\u00d755=\", newText: \"69\u00d778=\" },\n  { row: 0, col: 4, oldText: \"64\u00d770=\", newText: \"73\u00d711=\" },\n\n  { row: 4, col: 0, oldText: \"74\u00d795=\", newText: \"99\u00d713=\" },\n  { row: 4, col: 1, oldText: \"34\u00d765=\", newText: \"86\u00d724=\" },\n  { row: 4, col: 2, oldText: \"94\u00d750=\", newText: \"79\u00d766=\" },\n  { row: 4, col: 3, oldText: \"45\u00d716=\", newText: \"97\u00d739=\" },\n  { row: 4, col: 4, oldText: \"39\u00d721=\", newText: \"74\u00d724=\" },\n\n  { row: 9, col: 0, oldText: \"30\u00d745=\", newText: \"31\u00d719=\" },\n  { row: 9, col: 1, oldText: \"87\u00d758=\", newText: \"13\u00d758=\" },\n  { row: 9, col: 2, oldText: \"94\u00d750=\", newText: \"18\u00d795=\" },\n  { row: 9, col: 3, oldText: \"66\u00d714=\", newText: \"76\u00d773=\" },\n  { row: 9, col: 4, oldText: \"72\u00d780=\", newText: \"63\u00d757=\" },\n\n  { row: 14, col: 0, oldText: \"43\u00d757=\", newText: \"83\u00d742=\" },\n  { row: 14, col: 1, oldText: \"43\u00d740=\", newText: \"66\u00d770=\" },\n  { row: 14, col: 2, oldText: \"74\u00d795=\", newText: \"97\u00d744=\" },\n  { row: 14, col: 3, oldText: \"12\u00d739=\", newText: \"81\u00d717=\" },\n  { row: 14, col: 4, oldText: \"15\u00d796=\", newText: \"70\u00d771=\" },\n\n  { row: 19, col: 0, oldText: \"85\u00d733=\", newText: \"30\u00d781=\" },\n  { row: 19, col: 1, oldText: \"11\u00d790=\", newText: \"36\u00d747=\" },\n  { row: 19, col: 2, oldText: \"89\u00d756=\", newText: \"50\u00d723=\" },\n  { row: 19, col: 3, oldText: \"30\u00d782=\", newText: \"72\u00d728=\" },\n  { row: 19, col: 4, oldText: \"15\u00d718=\", newText: \"94\u00d793=\" },\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Scope a search to each target cell individually and replace the single\n// match in place, so formatting inherited from the original run survives.\nconst searchResultsByEdit = edits.map((edit) => {\n  const cell = table.getCell(edit.row, edit.col);\n  const results = cell.body.search(edit.oldText, { matchCase: true });\n  results.load(\"items\");\n  return results;\n});\n\nawait context.sync();\n\nsearchResultsByEdit.forEach((results, i) => {\n  const edit = edits[i];\n  if (results.items.length !== 1) {\n    throw new Error(\n      `Expected exactly one match for \"${edit.oldText}\" in cell (${edit.row},${edit.col}), found ${results.items.length}`\n    );\n  }\n  results.items[0].insertText(edit.newText, Word.InsertLocation.replace);\n});\n\nawait context.sync();\n", "ps1": "# Update the two-digit multiplication prompts in the worksheet table.\n# Each entry in $cells identifies a table cell by its 1-based (row, column)\n# position and gives the expected current text plus the replacement text.\n# Addressing cells by position (rather than a blind text search-and-replace)\n# is required because several prompts share the same original text\n# (e.g. \"94\u00d750=\" and \"74\u00d795=\" each occur twice) but must become different\n# values depending on which cell they are in.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$cells = @(\n    @{ Row = 1;  Col = 1; Old = \"58\u00d752=\"; New = \"18\u00d779=\" },\n    @{ Row = 1;  Col = 2; Old = \"99\u00d787=\"; New = \"46\u00d744=\" },\n    @{ Row = 1;  Col = 3; Old = \"77\u00d781=\"; New = \"72\u00d723=\" },\n    @{ Row = 1;  Col = 4; Old = \"98\u00d755=\"; New = \"69\u00d778=\" },\n    @{ Row = 1;  Col = 5; Old = \"64\u00d770=\"; New = \"73\u00d711=\" },\n\n    @{ Row = 5;  Col = 1; Old = \"74\u00d795=\"; New = \"99\u00d713=\" },\n    @{ Row = 5;  Col = 2; Old = \"34\u00d765=\"; New = \"86\u00d724=\" },\n    @{ Row = 5;  Col = 3; Old = \"94\u00d750=\"; New = \"79\u00d766=\" },\n    @{ Row = 5;  Col = 4; Old = \"45\u00d716=\"; New = \"97\u00d739=\" },\n    @{ Row = 5;  Col = 5; Old = \"39\u00d721=\"; New = \"74\u00d724=\" },\n\n    @{ Row = 10; Col = 1; Old = \"30\u00d745=\"; New = \"31\u00d719=\" },\n    @{ Row = 10; Col = 2; Old = \"87\u00d758=\"; New = \"13\u00d758=\" },\n    @{ Row = 10; Col = 3; Old = \"94\u00d750=\"; New = \"18\u00d795=\" },\n    @{ Row = 10; Col = 4; Old = \"66\u00d714=\"; New = \"76\u00d773=\" },\n    @{ Row = 10; Col = 5; Old = \"72\u00d780=\"; New = \"63\u00d757=\" },\n\n    @{ Row = 15; Col = 1; Old = \"43\u00d757=\"; New = \"83\u00d742=\" },\n    @{ Row = 15; Col = 2; Old = \"43\u00d740=\"; New = \"66\u00d770=\" },\n    @{ Row = 15; Col = 3; Old = \"74\u00d795=\"; New = \"97\u00d744=\" },\n    @{ Row = 15; Col = 4; Old = \"12\u00d739=\"; New = \"81\u00d717=\" },\n    @{ Row = 15; Col = 5; Old = \"15\u00d796=\"; New = \"70\u00d771=\" },\n\n    @{ Row = 20; Col = 1; Old = \"85\u00d733=\"; New = \"30\u00d781=\" },\n    @{ Row = 20; Col = 2; Old = \"11\u00d790=\"; New = \"36\u00d747=\" },\n    @{ Row = 20; Col = 3; Old = \"89\u00d756=\"; New = \"50\u00d723=\" },\n    @{ Row = 20; Col = 4; Old = \"30\u00d782=\"; New = \"72\u00d728=\" },\n    @{ Row = 20; Col = 5; Old = \"15\u00d718=\"; New = \"94\u00d793=\" }\n)\n\nforeach ($cell in $cells) {\n    $rng = $t.Cell($cell.Row, $cell.Col).Range\n    # Trim the trailing cell-mark characters so the comparison / assignment\n    # only touches the visible text, then write the new value in place.\n    $cellText = $rng.Text.TrimEnd([char]7, [char]13)\n    if ($cellText -ne $cell.Old) {\n        Write-Output (\"Warning: cell ({0},{1}) expected '{2}' but found '{3}'\" -f $cell.Row, $cell.Col, $cell.Old, $cellText)\n    }\n    $rng.Text = $cell.New\n}\n"}
